$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3:D18").ClearContents()

$ws.Range("D3:D18").Select()
